# "Add files via upload" - extend the mobility dataset on sheet "mobility"
# with 13 more days of data (2020-10-19 .. 2020-10-31), in columns KA:KM,
# for each of the 5 existing data rows. This grows the used range from
# A1:JZ5 to A1:KM5; Excel updates the sheet dimension/row spans on its own.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from JZ1 (the last existing date header) onto
# the new header cells KA1:KM1 so the new date serials inherit the same
# number format (style index "1" = date format) instead of General.
$ws.Range("JZ1").Copy()
$ws.Range("KA1:KM1").PasteSpecial(-4122)

# Row 1
$ws.Range("KA1").Value = 44123
$ws.Range("KB1").Value = 44124
$ws.Range("KC1").Value = 44125
$ws.Range("KD1").Value = 44126
$ws.Range("KE1").Value = 44127
$ws.Range("KF1").Value = 44128
$ws.Range("KG1").Value = 44129
$ws.Range("KH1").Value = 44130
$ws.Range("KI1").Value = 44131
$ws.Range("KJ1").Value = 44132
$ws.Range("KK1").Value = 44133
$ws.Range("KL1").Value = 44134
$ws.Range("KM1").Value = 44135

# Row 2
$ws.Range("KA2").Value = 50.52
$ws.Range("KB2").Value = 50.32
$ws.Range("KC2").Value = 51.62
$ws.Range("KD2").Value = 53.26
$ws.Range("KE2").Value = 62.07
$ws.Range("KF2").Value = 64.22
$ws.Range("KG2").Value = 47.94
$ws.Range("KH2").Value = 48.25
$ws.Range("KI2").Value = 50.23
$ws.Range("KJ2").Value = 49.65
$ws.Range("KK2").Value = 46.44
$ws.Range("KL2").Value = 55.89
$ws.Range("KM2").Value = 62.81

# Row 3
$ws.Range("KA3").Value = 36.25
$ws.Range("KB3").Value = 37.58
$ws.Range("KC3").Value = 36.15
$ws.Range("KD3").Value = 36.49
$ws.Range("KE3").Value = 41.28
$ws.Range("KF3").Value = 37.78
$ws.Range("KG3").Value = 34.3
$ws.Range("KH3").Value = 36.39
$ws.Range("KI3").Value = 33.69
$ws.Range("KJ3").Value = 36.95
$ws.Range("KK3").Value = 36.46
$ws.Range("KL3").Value = 42.93
$ws.Range("KM3").Value = 37.49

# Row 4
$ws.Range("KA4").Value = 61.54
$ws.Range("KB4").Value = 66.04
$ws.Range("KC4").Value = 63.83
$ws.Range("KD4").Value = 61.72
$ws.Range("KE4").Value = 65.99
$ws.Range("KF4").Value = 56.97
$ws.Range("KG4").Value = 48.29
$ws.Range("KH4").Value = 68.14
$ws.Range("KI4").Value = 62.74
$ws.Range("KJ4").Value = 65.06
$ws.Range("KK4").Value = 54.05
$ws.Range("KL4").Value = 64.17
$ws.Range("KM4").Value = 58.16

# Row 5
$ws.Range("KA5").Value = 63.8
$ws.Range("KB5").Value = 61.37
$ws.Range("KC5").Value = 64.18
$ws.Range("KD5").Value = 61.66
$ws.Range("KE5").Value = 66.86
$ws.Range("KF5").Value = 63.98
$ws.Range("KG5").Value = 51.12
$ws.Range("KH5").Value = 60.42
$ws.Range("KI5").Value = 60.4
$ws.Range("KJ5").Value = 58.89
$ws.Range("KK5").Value = 64.03
$ws.Range("KL5").Value = 69.11
$ws.Range("KM5").Value = 64.31
